$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '30.055.00'
$ws.Range("E2").Value = '  -0.55%  '

# Row 3
$ws.Range("D3").Value = '1.916.80'
$ws.Range("E3").Value = '  +0.52%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.50'
$ws.Range("D5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.02%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5042'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.17%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4022'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.00%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08296'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.02%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.107'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.20%  '

# Row 11
$ws.Range("E11").Value = '  -1.43%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.15'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.95%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.426'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.76%  '

# Row 14
$ws.Range("D14").Value = '1.913.94'
$ws.Range("E14").Value = '  +0.33%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.293'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.95%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9989'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.21%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.09'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.95%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001098'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.38%  '

# Row 19
$ws.Range("E19").Value = '  -2.43%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.86%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9999'
$ws.Range("D21").Style = "Normal"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.958'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.72%  '

# Row 23
$ws.Range("D23").Value = '30.099.80'
$ws.Range("E23").Value = '  -0.42%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.31%  '

# Row 25
$ws.Range("E25").Value = '  -0.15%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '22.35'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.42%  '

# Row 27
$ws.Range("D27").Value = '2.134.28'
$ws.Range("E27").Value = '  +0.40%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.24%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.283'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.41%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '129.14'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.50%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.126'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.36%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1039'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.83%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.997'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.84%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.797'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.88%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02449'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.82%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.358'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.01%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06415'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.35%  '

# Row 38
$ws.Range("B38").Value = 'TheSandbox'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6660'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.06%  '

# Row 39
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2162'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.83%  '

# Row 40
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.195'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.72%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.706'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.97%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.45'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.94%  '

# Row 43
$ws.Range("E43").Value = '  -0.78%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.210'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.18%  '

# Row 45
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6089'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.57%  '

# Row 46
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.31'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.89%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.641'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.13%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '123.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.24%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.214'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.53%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '78.70'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.70%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.127'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.68%  '
